$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "SUPERINTENDÊNCIA"
$ws.Range("H4").Value = "SUPERINTENDÊNCIA"
$ws.Range("G9").Value = "SUPORTE OPERACIONAL"
$ws.Range("H9").Value = "PROTOCOLO"
$ws.Range("G10").Value = "EXTERNO"
$ws.Range("H10").Value = "AUDITORIA MUNICIPAL"
$ws.Range("G11").Value = "EXTERNO"
$ws.Range("H11").Value = "LANCHONETE"
$ws.Range("G12").Value = "EXTERNO"
$ws.Range("H12").Value = "JURÍDICO"
$ws.Range("G13").Value = "SETOR DE CONTRATUALIZAÇÃO E REGULAÇÃO"
$ws.Range("H13").Value = "SETOR DE CONTRATUALIZAÇÃO E REGULAÇÃO"
$ws.Range("H53").Value = "SUPERINTENDÊNCIA"
$ws.Range("H54").Value = "SUPERINTENDÊNCIA"
$ws.Range("H55").Value = "SUPERINTENDÊNCIA"
$ws.Range("H56").Value = "SUPERINTENDÊNCIA"
$ws.Range("H57").Value = "SUPERINTENDÊNCIA"
$ws.Range("H58").Value = "SUPERINTENDÊNCIA"
$ws.Range("H59").Value = "SUPERINTENDÊNCIA"
$ws.Range("H60").Value = "SUPERINTENDÊNCIA"
$ws.Range("H61").Value = "SUPERINTENDÊNCIA"
$ws.Range("H62").Value = "SUPERINTENDÊNCIA"
$ws.Range("H63").Value = "SUPERINTENDÊNCIA"
$ws.Range("H64").Value = "SUPERINTENDÊNCIA"
$ws.Range("H65").Value = "SUPERINTENDÊNCIA"
$ws.Range("G66").Value = "FATURAMENTO HOSPITALAR"
$ws.Range("H66").Value = "FATURAMENTO HOSPITALAR"
$ws.Range("C70").Value = "Chefia UDIS"
$ws.Range("H110").Value = "UNIDADE DE NUTRIÇÃO CLÍNICA"
$ws.Range("H111").Value = "UNIDADE DE NUTRIÇÃO CLÍNICA"
$ws.Range("H112").Value = "UNIDADE DE NUTRIÇÃO CLÍNICA"
$ws.Range("H113").Value = "UNIDADE DE NUTRIÇÃO CLÍNICA"
$ws.Range("H114").Value = "UNIDADE DE NUTRIÇÃO CLÍNICA"
$ws.Range("H127").Value = "MORGE"
$ws.Range("C179").Value = " UMC-T-Recepcao-Maternidade"
$ws.Range("H202").Value = "UNIDADE DE APOIO Á GESTÃO EM INFERMAGEM"
$ws.Range("H203").Value = "UNIDADE DE APOIO Á GESTÃO EM INFERMAGEM"
$ws.Range("H204").Value = "UNIDADE DE APOIO Á GESTÃO EM INFERMAGEM"
$ws.Range("H205").Value = "UNIDADE DE APOIO Á GESTÃO EM INFERMAGEM"
$ws.Range("H241").Value = "UNIDADE DE ALMOXARIFADO E CONTROLE DE ESTOQUES"
$ws.Range("H242").Value = "UNIDADE DE ALMOXARIFADO E CONTROLE DE ESTOQUES"
$ws.Range("H243").Value = "UNIDADE DE ALMOXARIFADO E CONTROLE DE ESTOQUES"
$ws.Range("H244").Value = "UNIDADE DE ALMOXARIFADO E CONTROLE DE ESTOQUES"
$ws.Range("G304").Value = "SUPORTE OPERACIONAL"
$ws.Range("H304").Value = "SUPORTE OPERACIONAL"
$ws.Range("G305").Value = "SUPORTE OPERACIONAL"
$ws.Range("H305").Value = "SUPORTE OPERACIONAL"
$ws.Range("G306").Value = "SUPORTE OPERACIONAL"
$ws.Range("H306").Value = "SUPORTE OPERACIONAL"
$ws.Range("G307").Value = "SUPORTE OPERACIONAL"
$ws.Range("H307").Value = "GUARITA"
$ws.Range("G308").Value = "SUPORTE OPERACIONAL"
$ws.Range("H308").Value = "TRANSPORTES"
$ws.Range("G309").Value = "SUPORTE OPERACIONAL"
$ws.Range("H309").Value = "RECEPÇÕES"
$ws.Range("G310").Value = "SUPORTE OPERACIONAL"
$ws.Range("H310").Value = "RECEPÇÕES"
$ws.Range("G311").Value = "SUPORTE OPERACIONAL"
$ws.Range("H311").Value = "RECEPÇÕES"
$ws.Range("G312").Value = "SUPORTE OPERACIONAL"
$ws.Range("H312").Value = "RECEPÇÕES"
$ws.Range("G313").Value = "SUPORTE OPERACIONAL"
$ws.Range("H313").Value = "RECEPÇÕES"
$ws.Range("G314").Value = "SUPORTE OPERACIONAL"
$ws.Range("H314").Value = "RECEPÇÕES"
$ws.Range("G315").Value = "SUPORTE OPERACIONAL"
$ws.Range("H315").Value = "RECEPÇÕES"
$ws.Range("G316").Value = "SUPORTE OPERACIONAL"
$ws.Range("H316").Value = "RECEPÇÕES"
$ws.Range("G317").Value = "SUPORTE OPERACIONAL"
$ws.Range("H317").Value = "RECEPÇÕES"
$ws.Range("G318").Value = "SUPORTE OPERACIONAL"
$ws.Range("H318").Value = "RECEPÇÕES"
$ws.Range("G319").Value = "SUPORTE OPERACIONAL"
$ws.Range("H319").Value = "RECEPÇÕES"
$ws.Range("G320").Value = "SUPORTE OPERACIONAL"
$ws.Range("H320").Value = "RECEPÇÕES"
$ws.Range("G321").Value = "SUPORTE OPERACIONAL"
$ws.Range("H321").Value = "RECEPÇÕES"
$ws.Range("G322").Value = "SUPORTE OPERACIONAL"
$ws.Range("H322").Value = "RECEPÇÕES"
$ws.Range("G323").Value = "SUPORTE OPERACIONAL"
$ws.Range("H323").Value = "RECEPÇÕES"
$ws.Range("E341").Value = "SUPERINTENDÊNCIA"
$ws.Range("K341").Value = "P"
$ws.Range("E350").Value = "SUPERINTENDÊNCIA"
$ws.Range("C353").Value = "UMC-1-Teste da orelhinha-U01-060"
